$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.576331333333333
$ws.Range("H2").Value = 4.728994
$ws.Range("I2").Value = 0.01463337290888519
$ws.Range("J2").Value = 0.01463337290888519
$ws.Range("M2").Value = 0.09920366666666665
$ws.Range("N2").Value = 0.297611
$ws.Range("O2").Value = 0.641640866873065
$ws.Range("P2").Value = 0.6416408668730651
$ws.Range("Q2").Value = 0.1563778481482222
$ws.Range("R2").Value = 1.407400633334
$ws.Range("S2").Value = 0.009389370078533917
$ws.Range("T2").Value = 0.00938937007853392
$ws.Range("G3").Value = 1.576331333333333
$ws.Range("H3").Value = 4.728994
$ws.Range("I3").Value = 0.01463337290888519
$ws.Range("J3").Value = 0.01463337290888519
$ws.Range("O3").Value = 0.178081099028088
$ws.Range("P3").Value = 0.1780810990280881
$ws.Range("Q3").Value = 0.04340113060066667
$ws.Range("R3").Value = 0.390610175406
$ws.Range("S3").Value = 0.002605927130102124
$ws.Range("T3").Value = 0.002605927130102124
$ws.Range("G4").Value = 1.576331333333333
$ws.Range("H4").Value = 4.728994
$ws.Range("I4").Value = 0.01463337290888519
$ws.Range("J4").Value = 0.01463337290888519
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.02787266666666667
$ws.Range("N4").Value = 0.083618
$ws.Range("O4").Value = 0.180278034098847
$ws.Range("P4").Value = 0.180278034098847
$ws.Range("Q4").Value = 0.04393655781022222
$ws.Range("R4").Value = 0.395429020292
$ws.Range("S4").Value = 0.002638075700249148
$ws.Range("T4").Value = 0.002638075700249148
$ws.Range("G5").Value = 74.31489566666666
$ws.Range("I5").Value = 0.6898788078237544
$ws.Range("J5").Value = 0.6898788078237544
$ws.Range("M5").Value = 0.09920366666666665
$ws.Range("N5").Value = 0.297611
$ws.Range("O5").Value = 0.641640866873065
$ws.Range("P5").Value = 0.6416408668730651
$ws.Range("Q5").Value = 7.372310138084109
$ws.Range("R5").Value = 66.35079124275698
$ws.Range("S5").Value = 0.4426544362893904
$ws.Range("T5").Value = 0.4426544362893905
$ws.Range("G6").Value = 74.31489566666666
$ws.Range("I6").Value = 0.6898788078237544
$ws.Range("J6").Value = 0.6898788078237544
$ws.Range("O6").Value = 0.178081099028088
$ws.Range("P6").Value = 0.1780810990280881
$ws.Range("Q6").Value = 2.046112022390333
$ws.Range("S6").Value = 0.1228543762934413
$ws.Range("T6").Value = 0.1228543762934413
$ws.Range("G7").Value = 74.31489566666666
$ws.Range("I7").Value = 0.6898788078237544
$ws.Range("J7").Value = 0.6898788078237544
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.02787266666666667
$ws.Range("N7").Value = 0.083618
$ws.Range("O7").Value = 0.180278034098847
$ws.Range("P7").Value = 0.180278034098847
$ws.Range("Q7").Value = 2.071354315285111
$ws.Range("R7").Value = 18.642188837566
$ws.Range("S7").Value = 0.1243699952409227
$ws.Range("T7").Value = 0.1243699952409227
$ws.Range("G8").Value = 31.83044066666666
$ws.Range("H8").Value = 95.491322
$ws.Range("I8").Value = 0.2954878192673605
$ws.Range("J8").Value = 0.2954878192673605
$ws.Range("M8").Value = 0.09920366666666665
$ws.Range("N8").Value = 0.297611
$ws.Range("O8").Value = 0.641640866873065
$ws.Range("P8").Value = 0.6416408668730651
$ws.Range("Q8").Value = 3.15769642574911
$ws.Range("R8").Value = 28.41926783174199
$ws.Range("S8").Value = 0.1895970605051407
$ws.Range("T8").Value = 0.1895970605051408
$ws.Range("G9").Value = 31.83044066666666
$ws.Range("H9").Value = 95.491322
$ws.Range("I9").Value = 0.2954878192673605
$ws.Range("J9").Value = 0.2954878192673605
$ws.Range("O9").Value = 0.178081099028088
$ws.Range("P9").Value = 0.1780810990280881
$ws.Range("Q9").Value = 0.8763875228753333
$ws.Range("R9").Value = 7.887487705878001
$ws.Range("S9").Value = 0.0526207956045446
$ws.Range("T9").Value = 0.05262079560454461
$ws.Range("G10").Value = 31.83044066666666
$ws.Range("H10").Value = 95.491322
$ws.Range("I10").Value = 0.2954878192673605
$ws.Range("J10").Value = 0.2954878192673605
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.02787266666666667
$ws.Range("N10").Value = 0.083618
$ws.Range("O10").Value = 0.180278034098847
$ws.Range("P10").Value = 0.180278034098847
$ws.Range("Q10").Value = 0.8871992625551111
$ws.Range("R10").Value = 7.984793362995999
$ws.Range("S10").Value = 0.05326996315767516
$ws.Range("T10").Value = 0.05326996315767516
